# Add "hydrogen combined cycle" as a new power plant type (#99)
#
# - Rename the existing "hydrogen" row (A24) to "hydrogen combustion turbine"
#   and give it the black, vertically-centered font used for the new rows.
# - Add a new row (A25/B25) "hydrogen combined cycle" whose value mirrors
#   "natural gas combined cycle" (B4) via a formula, using the same style.
# - FoSYCRpUNL becomes the active sheet/selected tab (it was "About" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FoSYCRpUNL")

# Rename the "hydrogen" entry to "hydrogen combustion turbine".
$ws.Range("A24").Value = "hydrogen combustion turbine"

# Apply the new black / vertically-centered font to the relabeled cell.
$ws.Range("A24").Font.Color = 0
$ws.Range("A24").VerticalAlignment = -4108

# Add the new "hydrogen combined cycle" row, value = same as natural gas
# combined cycle (B4), wired as a formula so it tracks that input.
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25").Formula = "=B4"

# Clone A24's formatting onto A25 via copy/paste-special so both rows end
# up sharing a single cell style (rather than each cell accumulating its
# own separate style record).
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The FoSYCRpUNL sheet becomes the active tab/selected cell in this revision.
$ws.Activate()
$ws.Range("B26").Select()
